$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Mark existing experiments (rows 29-34) as "Done!" in the status column G
# ---------------------------------------------------------------------------
$ws.Range("G29:G34").Value = "Done!"

# ---------------------------------------------------------------------------
# Helper: write a rich-text "Experiment description" cell for a new
# regularization-weight experiment row. The text alternates bold / normal
# runs, mirroring the pattern used by the existing EXP29-EXP32 rows.
# ---------------------------------------------------------------------------
function Set-ExperimentDescription {
    param($Cell, [string]$Weight)

    $prefix = "Just TPR no LSTM in `nphrase embedding layer `nbatchsize = 40. "
    $text = $prefix + "With visualizations. With regularization. Regularization weights=$Weight [Resuming from latest successful commit, running from QA_TPR_for_Run_TPRregularizationFinal]. "

    $Cell.Value = $text

    $pos = 1
    $len1 = $prefix.Length
    $Cell.Characters($pos, $len1).Font.Name = "Calibri"
    $Cell.Characters($pos, $len1).Font.Size = 11
    $pos += $len1

    $len2 = "With visualizations".Length
    $Cell.Characters($pos, $len2).Font.Bold = $true
    $pos += $len2

    $len3 = 2
    $Cell.Characters($pos, $len3).Font.Name = "Calibri"
    $Cell.Characters($pos, $len3).Font.Size = 11
    $pos += $len3

    $len4 = "With regularization".Length
    $Cell.Characters($pos, $len4).Font.Bold = $true
    $pos += $len4

    $len5 = 2
    $Cell.Characters($pos, $len5).Font.Name = "Calibri"
    $Cell.Characters($pos, $len5).Font.Size = 11
    $pos += $len5

    $weightRun = "Regularization weights=$Weight"
    $len6 = $weightRun.Length
    $Cell.Characters($pos, $len6).Font.Bold = $true
    $pos += $len6

    $tail = " [Resuming from latest successful commit, running from QA_TPR_for_Run_TPRregularizationFinal]. "
    $len7 = $tail.Length
    $Cell.Characters($pos, $len7).Font.Name = "Calibri"
    $Cell.Characters($pos, $len7).Font.Size = 11
}

# ---------------------------------------------------------------------------
# 2. Append three new experiment rows (35, 36, 37) for EXP33, EXP34, EXP35 —
#    regularized-TPR runs with regularization weight 0.0001, 0.001, 0.01.
#    Start by copying the formatting of row 34 (fill/border/wrap/date style)
#    into the new rows (columns A-F hold the experiment data, K/L are the
#    still-empty "comments" / "date" columns; G-J are intentionally left
#    untouched since these new runs have not been marked done yet).
# ---------------------------------------------------------------------------
$ws.Range("A34:F34").Copy($ws.Range("A35:F35"))
$ws.Range("K34:L34").Copy($ws.Range("K35:L35"))
$ws.Range("A34:F34").Copy($ws.Range("A36:F36"))
$ws.Range("K34:L34").Copy($ws.Range("K36:L36"))
$ws.Range("A34:F34").Copy($ws.Range("A37:F37"))
$ws.Range("K34:L34").Copy($ws.Range("K37:L37"))

$ws.Rows.Item(35).RowHeight = 180
$ws.Rows.Item(36).RowHeight = 180
$ws.Rows.Item(37).RowHeight = 165

# -- Row 35: EXP33, cF/cR = 0.0001, run_id 30, pane 2, machine "DLDGX / 1"
Set-ExperimentDescription $ws.Range("A35") "0.0001"
$ws.Range("B35").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --justTPR True --TPRregularizer1 True --TPRvis True --cF 0.0001 --cR 0.0001 --batch_size 40 --run_id 30 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP33.txt"
$ws.Range("C35").Value = "DLDGX / 1"
$ws.Range("D35").Value = "EXP33.txt"
$ws.Range("E35").Value = 30
$ws.Range("F35").Value = 2

# -- Row 36: EXP34, cF/cR = 0.001, run_id 31, pane 3, machine "DLDGX / 2"
Set-ExperimentDescription $ws.Range("A36") "0.001"
$ws.Range("B36").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --justTPR True --TPRregularizer1 True --TPRvis True --cF 0.001 --cR 0.001 --batch_size 40 --run_id 31 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP34.txt"
$ws.Range("C36").Value = "DLDGX / 2"
$ws.Range("D36").Value = "EXP34.txt"
$ws.Range("E36").Value = 31
$ws.Range("F36").Value = 3

# -- Row 37: EXP35, cF/cR = 0.01, run_id 32, pane 4, machine "DLDGX / 3"
Set-ExperimentDescription $ws.Range("A37") "0.01"
$ws.Range("B37").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --justTPR True --TPRregularizer1 True --TPRvis True --cF 0.01 --cR 0.01 --batch_size 40 --run_id 32 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP35.txt"
$ws.Range("C37").Value = "DLDGX / 3"
$ws.Range("D37").Value = "EXP35.txt"
$ws.Range("E37").Value = 32
$ws.Range("F37").Value = 4

# ---------------------------------------------------------------------------
# 3. Move the frozen-pane view / selection down to the newly-added rows.
# ---------------------------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("G34").Select() | Out-Null

Write-Output "edit complete"
